# Moved to U.S. version 2.1.1 as baseline
#
# The "Notes" explanation on the About sheet is condensed from a
# five-line passage down to a single sentence, and as a consequence the
# BVTStL sheet's boolean flags are updated: only aircraft remain exempt
# from the LCFS, so rail and ships are now marked as subject to it (1).

$wb = $excel.ActiveWorkbook

# ---- "About" sheet -------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Replace the old 5-line note (rows 15-19) with the new single sentence
# and clear out what used to be its continuation lines.
$about.Range("A15").Value = "Based on the California LCFS, we choose to exempt aircraft."
$about.Range("A16").ClearContents()
$about.Range("A17").ClearContents()
$about.Range("A18").ClearContents()
$about.Range("A19").ClearContents()

# ---- "BVTStL" sheet --------------------------------------------------
$bvtstl = $wb.Worksheets.Item("BVTStL")

# Rail and ships are no longer exempt -- they're now subject to the LCFS.
$bvtstl.Range("B5").Value = 1
$bvtstl.Range("C5").Value = 1
$bvtstl.Range("B6").Value = 1
$bvtstl.Range("C6").Value = 1

[void]$bvtstl.Range("C4").Select()

# Leave "About" as the active/selected sheet (it was the sheet shown
# when the workbook was originally saved).
[void]$about.Activate()
[void]$about.Range("C18").Select()
